$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Day 3")
$ws4 = $wb.Worksheets.Item("Day 4")

# --- Day 3 ("T-ROW / Back" style day) -----------------------------------
# Values are written in this specific order so that new shared-strings
# are interned in the same sequence the original authors' edit produced.
$ws3.Range("E10").Value = "3x15"
$ws3.Range("E2").Value = "3x10"
$ws3.Range("E1").Value = "3x12"
$ws3.Range("E6").Value = "3x8e"
$ws3.Range("A10").Value = "Lat Pulldown (slow negative)"
$ws3.Range("A2").Value = "T-ROW Machine"
$ws3.Range("A1").Value = "Standing Cable Bicep Curls"
$ws3.Range("A6").Value = "Seated DB Hammer Curls"
$ws3.Range("A9").Value = "Kettlebell Squat-to-Upright Row"
$ws3.Range("A3").Value = "Seated Back Extensions"
$ws3.Range("A5").Value = "Bent over Barbell Rows"
$ws3.Range("A7").Value = "Isometric Rear Delt Fly"
$ws3.Range("E7").Value = "4x60s"
$ws3.Range("E5").Value = "3x10"
$ws3.Range("E9").Value = "3x10"

# --- Day 4 (circuit day) -------------------------------------------------
$ws4.Range("A1").Value = "Low-to-High Side Chops (Kettlebell) 3x8e"
$ws4.Range("A2").Value = "Weighted Box Squats 3x10"
$ws4.Range("A3").Value = "Plank Hold 3x1min."
$ws4.Range("A4").Value = "Medicine Ball Alt. Pushups 3x8e"
$ws4.Range("A5").Value = "Spiders (Knee to Elbow) 3x8e"
$ws4.Range("A6").Value = "Raised Leg Lifts 3x10"
$ws4.Range("B3").Value = "3x1min"
$ws4.Range("C1").Value = "*"
$ws4.Range("C3").Value = "(Bodyweight Plank)"
$ws4.Range("D6").Value = "(Raised Leg Lifts)"
$ws4.Range("D5").Value = "(Spiders, Knee to Elbow)"
$ws4.Range("D1").Value = "(Low-to-High Side Chops, Kettlebell)"
$ws4.Range("C2").Value = "(Dumbbell Squat, Assisted)"
$ws4.Range("C4").Value = "(Medicine Ball Pushup, Alternating)"
$ws4.Range("B1").Value = "3x8e"
$ws4.Range("B2").Value = "3x10"
$ws4.Range("B4").Value = "3x8e"
$ws4.Range("B5").Value = "3x8e"
$ws4.Range("B6").Value = "3x10"
$ws4.Range("C5").Value = "*"
$ws4.Range("C6").Value = "*"

# Column B on Day 4 (the rep-scheme column) is centered.
$ws4.Range("B1:B6").HorizontalAlignment = -4108

# Day 4 column A is widened to fit the longest exercise label.
$ws4.Columns.Item(1).ColumnWidth = 37.45

# --- Sheet views / selections --------------------------------------------
[void]$ws3.Range("B21").Select()
[void]$ws4.Range("C5").Select()

# --- Page setup on Day 2 --------------------------------------------------
$ws2 = $wb.Worksheets.Item("Day 2")
$ws2.PageSetup.Orientation = 1

# --- Active sheet moves from Day 2 to Day 4 -------------------------------
$ws4.Activate()
